$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...Notowidigdo) Conditionally Accepted at ..." ->
#    "...Notowidigdo) forthcoming at ..."
#    (the journal-status phrase "Conditionally Accepted" becomes "forthcoming")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Conditionally Accepted", $true, $false, $false, $false, $false, $true, 1, $false, "forthcoming", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert a new, empty BodyText paragraph right after the paragraph that
#    ends in "...128(3): 1123- 1167, August 2013." (and before the
#    "Does Online Search..." paragraph). Using "^p" in the replacement text
#    of a Find/Replace mints a clean new <w:p> that inherits the
#    surrounding BodyText paragraph formatting, with no stray runs.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("1123- 1167, August 2013.", $true, $false, $false, $false, $false, $true, 1, $false, "1123- 1167, August 2013.^p", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Drop the stale <w:lastRenderedPageBreak/> marker sitting in front of the
#    "PROFESSIONAL RESPONSIBILITIES" heading run. Re-asserting the run text
#    in place rebuilds the run and sheds the obsolete page-break bookmark.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("PROFESSIONAL RESPONSIBILITIES", $true, $false, $false, $false, $false, $true, 1, $false, "PROFESSIONAL RESPONSIBILITIES", 2) | Out-Null
